{"js": "// Replace each \"divided by\" expression in the practice table with its\n// updated counterpart. Every source/target string below is unique within\n// the document, so a direct search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"781\u00f78=\", \"524\u00f77=\"],\n  [\"541\u00f74=\", \"923\u00f77=\"],\n  [\"774\u00f79=\", \"333\u00f75=\"],\n  [\"987\u00f78=\", \"428\u00f74=\"],\n  [\"542\u00f72=\", \"398\u00f76=\"],\n  [\"291\u00f79=\", \"152\u00f75=\"],\n  [\"488\u00f76=\", \"766\u00f78=\"],\n  [\"474\u00f79=\", \"856\u00f76=\"],\n  [\"934\u00f79=\", \"884\u00f77=\"],\n  [\"604\u00f77=\", \"238\u00f78=\"],\n  [\"158\u00f73=\", \"588\u00f76=\"],\n  [\"850\u00f74=\", \"474\u00f75=\"],\n  [\"211\u00f79=\", \"153\u00f76=\"],\n  [\"576\u00f74=\", \"842\u00f78=\"],\n  [\"945\u00f78=\", \"680\u00f78=\"],\n  [\"842\u00f79=\", \"564\u00f79=\"],\n  [\"462\u00f74=\", \"951\u00f77=\"],\n  [\"267\u00f73=\", \"778\u00f78=\"],\n  [\"353\u00f77=\", \"846\u00f75=\"],\n  [\"486\u00f79=\", \"623\u00f78=\"],\n  [\"713\u00f74=\", \"861\u00f79=\"],\n  [\"465\u00f75=\", \"466\u00f73=\"],\n  [\"882\u00f72=\", \"413\u00f79=\"],\n  [\"663\u00f75=\", \"345\u00f77=\"],\n  [\"614\u00f73=\", \"736\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"divided by\" expression in the practice table with its\n# updated counterpart. Every source/target string is unique within the\n# document, so Find/Replace (wdReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"781\u00f78=\", \"524\u00f77=\"),\n    @(\"541\u00f74=\", \"923\u00f77=\"),\n    @(\"774\u00f79=\", \"333\u00f75=\"),\n    @(\"987\u00f78=\", \"428\u00f74=\"),\n    @(\"542\u00f72=\", \"398\u00f76=\"),\n    @(\"291\u00f79=\", \"152\u00f75=\"),\n    @(\"488\u00f76=\", \"766\u00f78=\"),\n    @(\"474\u00f79=\", \"856\u00f76=\"),\n    @(\"934\u00f79=\", \"884\u00f77=\"),\n    @(\"604\u00f77=\", \"238\u00f78=\"),\n    @(\"158\u00f73=\", \"588\u00f76=\"),\n    @(\"850\u00f74=\", \"474\u00f75=\"),\n    @(\"211\u00f79=\", \"153\u00f76=\"),\n    @(\"576\u00f74=\", \"842\u00f78=\"),\n    @(\"945\u00f78=\", \"680\u00f78=\"),\n    @(\"842\u00f79=\", \"564\u00f79=\"),\n    @(\"462\u00f74=\", \"951\u00f77=\"),\n    @(\"267\u00f73=\", \"778\u00f78=\"),\n    @(\"353\u00f77=\", \"846\u00f75=\"),\n    @(\"486\u00f79=\", \"623\u00f78=\"),\n    @(\"713\u00f74=\", \"861\u00f79=\"),\n    @(\"465\u00f75=\", \"466\u00f73=\"),\n    @(\"882\u00f72=\", \"413\u00f79=\"),\n    @(\"663\u00f75=\", \"345\u00f77=\"),\n    @(\"614\u00f73=\", \"736\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
